# Swap the worker records currently on rows 17 and 18:
#   Row 17: 1047482742 / LEONARDO JAVIER VERA DORIA / 1708
#   Row 18: 45757837   / LORENA BEATRIZ DONADO LOPEZ / 1712
# becomes:
#   Row 17: 45757837   / LORENA BEATRIZ DONADO LOPEZ / 1712
#   Row 18: 1047482742 / LEONARDO JAVIER VERA DORIA / 1708

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = "45757837"
$ws.Range("D17").Value = "LORENA BEATRIZ DONADO LOPEZ"
$ws.Range("E17").Value = "1712"

$ws.Range("C18").Value = "1047482742"
$ws.Range("D18").Value = "LEONARDO JAVIER VERA DORIA"
$ws.Range("E18").Value = "1708"
